$d = $word.ActiveDocument

$replacements = @(
    @("544×4=", "548×7="),
    @("476×4=", "904×9="),
    @("766×2=", "155×5="),
    @("559×3=", "417×5="),
    @("568×7=", "379×5="),
    @("855×5=", "952×3="),
    @("817×2=", "836×7="),
    @("835×3=", "917×8="),
    @("263×2=", "728×3="),
    @("754×6=", "135×5="),
    @("687×4=", "396×2="),
    @("119×2=", "963×3="),
    @("850×4=", "134×9="),
    @("413×3=", "614×6="),
    @("508×6=", "224×3="),
    @("897×2=", "675×6="),
    @("614×5=", "844×7="),
    @("288×2=", "212×2="),
    @("565×2=", "354×7="),
    @("929×8=", "794×5="),
    @("471×8=", "621×4="),
    @("954×3=", "543×9="),
    @("652×7=", "272×9="),
    @("812×3=", "678×5="),
    @("636×7=", "782×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
